$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 341, shifting existing rows 341-375 down to 342-376
$ws.Rows.Item(341).Insert()

# Populate the new row 341 with the full record (copy of row 342's static fields + new data)
$ws.Cells.Item(341, 1).Value = 7
$ws.Cells.Item(341, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(341, 3).Value = "Ñuble"
$ws.Cells.Item(341, 4).Value = 45166
$ws.Cells.Item(341, 4).NumberFormat = $ws.Cells.Item(342, 4).NumberFormat
$ws.Cells.Item(341, 5).Value = 16
$ws.Cells.Item(341, 6).Value = "Fruta"
$ws.Cells.Item(341, 7).Value = 100101
$ws.Cells.Item(341, 8).Value = "Berries"
$ws.Cells.Item(341, 9).Value = 100101007
$ws.Cells.Item(341, 10).Value = "Kiwi"
$ws.Cells.Item(341, 11).Value = "Hayward"
$ws.Cells.Item(341, 12).Value = "Primera"
$ws.Cells.Item(341, 13).Value = 80
$ws.Cells.Item(341, 14).Value = 14000
$ws.Cells.Item(341, 15).Value = 14000
$ws.Cells.Item(341, 16).Value = 14000
$ws.Cells.Item(341, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(341, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(341, 19).Value = 778
$ws.Cells.Item(341, 20).Value = 18
